# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 4: Qualcomm Atheros AR9580 Wireless Network Adapter - 10.1.10.5
$ws.Range("C4").Value = 902
$ws.Range("D4").Value = 94.59999999999999

# Row 5: Qualcomm Atheros AR9580 Wireless Network Adapter - 3.0.2.201
$ws.Range("C5").Value = 1063

# Row 6: Totals
$ws.Range("C6").Value = 1966
